$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 (PPP012 / Yogesh Kumar JG): Pending Task shortened to a single item.
$ws.Range("D13").Value = "update_linkedin_with-photo"

# Row 15 (PPP014 / Naveen Bromiyo A R): Pending Task shortened to two items.
$ws.Range("D15").Value = "create_wordpress_blog_and_7articles, update_linkedin_with-photo"

# Footer timestamp refreshed for the new report run.
$ws.Range("A30").Value = "Generated: 2023-09-04 10:47:15 AM"
